$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Phase 0: style 19 (bold TOTAL-row %Chg format) donor is N21 itself,
# which also changes value later -> stamp its style onto the other
# style-19 targets BEFORE we touch N21.
# ---------------------------------------------------------------
$ws.Range("N21").Copy()
$ws.Range("E21").PasteSpecial(-4122)
$ws.Range("H21").PasteSpecial(-4122)
$ws.Range("K21").PasteSpecial(-4122)
$ws.Range("L21").PasteSpecial(-4122)
$ws.Range("M21").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Phase 1: cells whose style is unchanged - plain value assignment
# ---------------------------------------------------------------
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -14.285714285714
$ws.Range("N16").Value = -70
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -72.727272727272
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -44.827586206896
$ws.Range("I17").Value = 4
$ws.Range("L17").Value = 33.333333333333
$ws.Range("N17").Value = -63.636363636363
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -62.5
$ws.Range("M18").Value = -66.666666666666
$ws.Range("N18").Value = -85.714285714285
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = -23.076923076923
$ws.Range("F20").Value = 4
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 3
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -22.222222222222
$ws.Range("F21").Value = 41
$ws.Range("G21").Value = 63
$ws.Range("H21").Value = -34.920634920634
$ws.Range("I21").Value = 18
$ws.Range("L21").Value = 80
$ws.Range("M21").Value = 20
$ws.Range("N21").Value = -70.491803278688
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -75
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -61.538461538461
$ws.Range("I23").Value = 2
$ws.Range("C24").Value = 9
$ws.Range("E24").Value = -10
$ws.Range("G24").Value = 40
$ws.Range("H24").Value = 20
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 24
$ws.Range("H25").Value = -20
$ws.Range("M25").Value = 0
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -100
$ws.Range("J40").Value = 72
$ws.Range("K40").Value = -66.666666666666
$ws.Range("L40").Value = -80.54054054054
$ws.Range("M40").Value = -90.625
$ws.Range("N40").Value = -91.272727272727
$ws.Range("J41").Value = 161
$ws.Range("K41").Value = -14.814814814814
$ws.Range("L41").Value = -20.689655172413
$ws.Range("M41").Value = -51.212121212121
$ws.Range("N41").Value = -55.524861878453

# ---------------------------------------------------------------
# Phase 2: cells converting NUMBER -> TEXT placeholder (style 14/17)
# Must set NumberFormat=text + value BEFORE pasting the real style,
# otherwise re-applying the style after makes Excel re-parse "0" as a number.
# ---------------------------------------------------------------
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("A14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("F27").NumberFormat = "@"
$ws.Range("F27").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C29").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Phase 3: cells converting TEXT placeholder -> NUMBER (style 15/16/18/19)
# Set the numeric value directly, then paste the exact target style.
# ---------------------------------------------------------------
$ws.Range("L15").Value = -100
$ws.Range("H28").Copy()
$ws.Range("L15").PasteSpecial(-4122)
$ws.Range("I16").Value = 3
$ws.Range("F15").Copy()
$ws.Range("I16").PasteSpecial(-4122)
$ws.Range("L16").Value = 200
$ws.Range("H28").Copy()
$ws.Range("L16").PasteSpecial(-4122)
$ws.Range("M16").Value = -25
$ws.Range("H28").Copy()
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range("J17").Value = 11
$ws.Range("F15").Copy()
$ws.Range("J17").PasteSpecial(-4122)
$ws.Range("K17").Value = -63.636363636363
$ws.Range("H28").Copy()
$ws.Range("K17").PasteSpecial(-4122)
$ws.Range("M17").Value = 33.333333333333
$ws.Range("H28").Copy()
$ws.Range("M17").PasteSpecial(-4122)
$ws.Range("D18").Value = 4
$ws.Range("F15").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").Value = -100
$ws.Range("H28").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("J18").Value = 4
$ws.Range("F15").Copy()
$ws.Range("J18").PasteSpecial(-4122)
$ws.Range("K18").Value = -50
$ws.Range("H28").Copy()
$ws.Range("K18").PasteSpecial(-4122)
$ws.Range("L18").Value = 0
$ws.Range("H28").Copy()
$ws.Range("L18").PasteSpecial(-4122)
$ws.Range("I19").Value = 6
$ws.Range("F15").Copy()
$ws.Range("I19").PasteSpecial(-4122)
$ws.Range("J19").Value = 3
$ws.Range("F15").Copy()
$ws.Range("J19").PasteSpecial(-4122)
$ws.Range("K19").Value = 100
$ws.Range("H28").Copy()
$ws.Range("K19").PasteSpecial(-4122)
$ws.Range("L19").Value = 500
$ws.Range("H28").Copy()
$ws.Range("L19").PasteSpecial(-4122)
$ws.Range("M19").Value = 500
$ws.Range("H28").Copy()
$ws.Range("M19").PasteSpecial(-4122)
$ws.Range("N19").Value = -45.454545454545
$ws.Range("H28").Copy()
$ws.Range("N19").PasteSpecial(-4122)
$ws.Range("L20").Value = 200
$ws.Range("H28").Copy()
$ws.Range("L20").PasteSpecial(-4122)
$ws.Range("M20").Value = 200
$ws.Range("H28").Copy()
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range("N20").Value = -80
$ws.Range("H28").Copy()
$ws.Range("N20").PasteSpecial(-4122)
$ws.Range("J21").Value = 18
$ws.Range("C43").Copy()
$ws.Range("J21").PasteSpecial(-4122)
$ws.Range("K21").Value = 0
$ws.Range("N21").Copy()
$ws.Range("K21").PasteSpecial(-4122)
$ws.Range("J23").Value = 4
$ws.Range("F15").Copy()
$ws.Range("J23").PasteSpecial(-4122)
$ws.Range("K23").Value = -50
$ws.Range("H28").Copy()
$ws.Range("K23").PasteSpecial(-4122)
$ws.Range("L23").Value = 100
$ws.Range("H28").Copy()
$ws.Range("L23").PasteSpecial(-4122)
$ws.Range("M23").Value = 100
$ws.Range("H28").Copy()
$ws.Range("M23").PasteSpecial(-4122)
$ws.Range("I24").Value = 9
$ws.Range("F15").Copy()
$ws.Range("I24").PasteSpecial(-4122)
$ws.Range("J24").Value = 10
$ws.Range("F15").Copy()
$ws.Range("J24").PasteSpecial(-4122)
$ws.Range("K24").Value = -10
$ws.Range("H28").Copy()
$ws.Range("K24").PasteSpecial(-4122)
$ws.Range("L24").Value = 12.5
$ws.Range("H28").Copy()
$ws.Range("L24").PasteSpecial(-4122)
$ws.Range("M24").Value = 125
$ws.Range("H28").Copy()
$ws.Range("M24").PasteSpecial(-4122)
$ws.Range("I25").Value = 9
$ws.Range("F15").Copy()
$ws.Range("I25").PasteSpecial(-4122)
$ws.Range("J25").Value = 6
$ws.Range("F15").Copy()
$ws.Range("J25").PasteSpecial(-4122)
$ws.Range("K25").Value = 50
$ws.Range("H28").Copy()
$ws.Range("K25").PasteSpecial(-4122)
$ws.Range("L25").Value = 12.5
$ws.Range("H28").Copy()
$ws.Range("L25").PasteSpecial(-4122)
$ws.Range("L26").Value = -100
$ws.Range("H28").Copy()
$ws.Range("L26").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("F15").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("H28").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("J27").Value = 1
$ws.Range("F15").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("K27").Value = -100
$ws.Range("H28").Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("L27").Value = -100
$ws.Range("H28").Copy()
$ws.Range("L27").PasteSpecial(-4122)
$ws.Range("N28").Value = -100
$ws.Range("H28").Copy()
$ws.Range("N28").PasteSpecial(-4122)
$ws.Range("N29").Value = -100
$ws.Range("H28").Copy()
$ws.Range("N29").PasteSpecial(-4122)
